$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24; this shifts rows 24..124 down to 25..125
# and pushes the dimension ref to A1:R125.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record's data.
$ws.Cells.Item(24, 1).Value = 10
$ws.Cells.Item(24, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(24, 3).Value = "La Araucanía"
$ws.Cells.Item(24, 4).Value = 45222
$ws.Cells.Item(24, 5).Value = 9
$ws.Cells.Item(24, 6).Value = 100112022
$ws.Cells.Item(24, 7).Value = "Arveja Verde"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 50
$ws.Cells.Item(24, 11).Value = 35000
$ws.Cells.Item(24, 12).Value = 35000
$ws.Cells.Item(24, 13).Value = 35000
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región del Maule"
$ws.Cells.Item(24, 16).Value = 1400
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
